# Generate Report for Handback
#
# The handback report workbook tracks localization status per target
# locale. Running a new handback generation:
#   - flips each locale's Status from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview + each locale sheet)
#   - refreshes the "Latest Handback DateTime" for each locale
#   - clears the stale "Error Detail" (no longer an error once in sync)
#   - widens the Status column / narrows the (now empty) Error Detail
#     column to fit their new content

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.144371396019366
$overview.Columns.Item(6).ColumnWidth = 29.144371396019366

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-09-01 10:56:59"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.144371396019366
$zhcn.Columns.Item(16).ColumnWidth = 12.913719540550566

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-09-01 10:57:14"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.144371396019366
$dede.Columns.Item(16).ColumnWidth = 12.913719540550566
